$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Header fix: "unnamed: 1_level_1" -> "total"
$ws.Range("B2").Value = "total"

# Shift the numeric block up: rows that used to be header/footnote-only
# labels with no data of their own (old row 5 "situacao do domicilio", old
# row 8 "grandes regioes e unidades da federacao", and old row 41 "fonte:
# ibge...") are dropped, and every data row below moves up to close the gap.
$data = @{
    4 = @(1.63, 5, 3.4, 2.53, 2.27, 2.47)
    5 = @(1.6, 4.96, 3.57, 2.6, 2.33, 2.53)
    6 = @(4.6, 10.03, 7.67, 7.56, 7.86, 10.5)
    7 = @(4.15, 8.65, 7.76, 6.06, 5.51, 6.16)
    8 = @(9.01, 28.66, 22.02, 15.88, 14.05, 11.4)
    9 = @(10.62, 25.13, 26.42, 15.09, 10.66, 16.33)
    10 = @(11.48, 26.72, 26.53, 14.7, 13.07, 16.25)
    11 = @(10.5, 26.92, 22.47, 16.89, 21.91, 18.8)
    12 = @(7.19, 11.37, 10.99, 10.24, 10.43, 13.01)
    13 = @(13.54, 45.39, 31.13, 14.24, 19.51, 18.34)
    14 = @(11.76, 32.12, 21.47, 18.27, 9.88, 17.61)
    15 = @(3.9, 8.29, 5.84, 4.61, 5.59, 5.24)
    16 = @(17.85, 37.56, 20.07, 15.69, 23.18, 19.27)
    17 = @(13.19, 23.12, 20.68, 16.95, 19.22, 12.93)
    18 = @(6.94, 13.89, 10.08, 9.67, 10.85, 11.26)
    19 = @(10.71, 26.64, 16.44, 13.83, 15.27, 20.48)
    20 = @(11.65, 30.44, 28.6, 21.52, 24.63, 21.45)
    21 = @(6.05, 15.74, 12.71, 9.86, 12.08, 12.87)
    22 = @(10.84, 22.62, 17.95, 18.14, 17.07, 28.17)
    23 = @(9.18, 19.96, 13.01, 15.63, 10.71, 19.44)
    24 = @(5.69, 9.52, 10.45, 8.7, 8.01, 9.03)
    25 = @(2.8, 9.46, 7.17, 5.12, 3.82, 4.21)
    26 = @(4.95, 14.96, 11.49, 8.47, 5.7, 6.28)
    27 = @(10.31, 26.81, 17.35, 15.9, 14.95, 15.29)
    28 = @(5.72, 22.98, 15.87, 10.55, 9.88, 8.68)
    29 = @(4.23, 16.55, 13.09, 8.14, 5.78, 6.41)
    30 = @(3.2, 12.35, 8.35, 5.69, 4.86, 4.47)
    31 = @(4.65, 16.47, 14.37, 8.43, 7.49, 7.14)
    32 = @(7.51, 40.75, 21.16, 13.34, 11.46, 9.4)
    33 = @(5.22, 17.61, 11.84, 9.3, 7.36, 7.29)
    34 = @(3.19, 14.63, 8.44, 5.33, 4.49, 5.3)
    35 = @(7.25, 35.2, 18.29, 10.27, 8.65, 12.71)
    36 = @(7.01, 28, 22.74, 12.14, 7.89, 12.91)
    37 = @(4.99, 23.27, 10.62, 8.41, 7.3, 9.26)
    38 = @(6.46, 31.22, 31.68, 11.16, 12.18, 8.88)
}

foreach ($r in ($data.Keys | Sort-Object {[int]$_})) {
    $vals = $data[$r]
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $ws.Cells.Item([int]$r, 2 + $i).Value = $vals[$i]
    }
}

# Rows 39-41 (old rows 43-45) no longer belong to the table; clear them
$ws.Range("A39:G41").Clear()
